# update the final test result
#
# Rows 13, 15 and 16 used to be marked "failed" / "failed on test case 2" and
# are now passing, so their value becomes "passed" and the red/orange
# highlight fill is removed (copy the plain "passed" formatting from C14).
#
# Rows 11 and 12 are still failing, but now carry a longer explanation, so
# the cell text grows to two lines and Wrap Text is turned on for them
# (the existing highlight fill is left as-is).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Rows that flipped from failing to passing ------------------------------
$ws.Range("C13").Value = "passed"
$ws.Range("C15").Value = "passed"
$ws.Range("C16").Value = "passed"

$ws.Range("C14").Copy()
$ws.Range("C13").PasteSpecial(-4122)  # xlPasteFormats - drop the old fail highlight
$ws.Range("C15").PasteSpecial(-4122)  # xlPasteFormats - drop the old fail highlight
$ws.Range("C16").PasteSpecial(-4122)  # xlPasteFormats - drop the old fail highlight

# --- Rows that are still failing, with an updated / longer message ---------
$ws.Range("C11").Value = "failed.`nNested query is not supported yet."
$ws.Range("C12").Value = "failed.`nNested query is not supported yet."
$ws.Range("C11:C12").WrapText = $true

# --- View state: scroll position & active selection -------------------------
$ws.Activate() | Out-Null
$excel.ActiveWindow.ScrollRow = 9
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("F12").Select() | Out-Null
